$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.757.44"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.01%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.630.62"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.13%  "

$ws.Range("E4").Value = "  -0.55%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.19%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.501"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.40%  "

$ws.Range("E7").Value = "  -0.60%  "

$ws.Range("E8").Value = "  -0.89%  "

$ws.Range("E9").Value = "  -0.74%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.53"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.21%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0791"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.95%  "

$ws.Range("E12").Value = "  +0.39%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.856.33"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.09%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.611.23"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.43%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.553"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.39%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₃0761"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.45%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.86"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.67%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.761.33"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.07%  "

$ws.Range("E19").Value = "  -0.58%  "

$ws.Range("E20").Value = "  +0.37%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "191.46"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.08%  "

$ws.Range("E22").Value = "  -0.25%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.03%  "

$ws.Range("E24").Value = "  +1.83%  "

$ws.Range("E25").Value = "  -0.59%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "142.46"
$ws.Range("D26").Style = "Normal"

$ws.Range("E27").Value = "  +2.92%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.85"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.33%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.48"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.50%  "

$ws.Range("E30").Value = "  -0.39%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0494"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.46%  "

$ws.Range("E32").Value = "  -0.10%  "

$ws.Range("E33").Value = "  -0.63%  "

$ws.Range("E34").Value = "  +0.23%  "

$ws.Range("E35").Value = "  -0.55%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.906"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.18%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.138.46"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.67%  "

$ws.Range("E38").Value = "  -2.82%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.543"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.65%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0156"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.24%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.60%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.53"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.62%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.69"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.08%  "

$ws.Range("E44").Value = "  -0.71%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.805"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.38%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.766.66"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.28%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.29"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.43%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0511"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.22%  "

$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.45"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.96%  "

$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.418"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.05%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.35"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.01%  "
